$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "황구조(荒鉤爪)"
$ws.Range("B2").Value = "황구조의 혼(荒鉤爪の魂)"
$ws.Range("C2").Value = "10"
$ws.Range("D2").Value = "속식+2(早食い+2) + 고급 귀마개(高級耳栓)"
$ws.Range("A3").Value = "거합(居合)"
$ws.Range("B3").Value = "거합술【력】(居合術【力】)"
$ws.Range("C3").Value = "10"
$ws.Range("D3").Value = "발도술【력】(抜刀術【力】) + 납도술(納刀術)"
$ws.Range("A4").Value = "노(怒)"
$ws.Range("B4").Value = "역린(逆鱗)"
$ws.Range("C4").Value = "10"
$ws.Range("D4").Value = "화사장력+2(火事場力+2) + 근성(根性)"
$ws.Range("A5").Value = "암천(岩穿)"
$ws.Range("B5").Value = "암천의 혼(岩穿の魂)"
$ws.Range("C5").Value = "10"
$ws.Range("D5").Value = "화사장력+2(火事場力+2) + 주워먹기(拾い食い)"
$ws.Range("A6").Value = "포말(泡沫)"
$ws.Range("B6").Value = "포말의 춤(泡沫の舞)"
$ws.Range("C6").Value = "10"
$ws.Range("D6").Value = "거품상태【대】가 되지 않는다.`n회피행동을 반복(통상회피 3회, 스텝회피4회)하면 거품상태【소】가 되어 회피시의 무적시간이 늘어나고 스태미너 소비가 감소한다.`n(회피성능+1, 체술+1)"
$ws.Range("A7").Value = "역회심(裏会心)"
$ws.Range("B7").Value = "통한회심(痛恨会心)"
$ws.Range("C7").Value = "10"
$ws.Range("D7").Value = "마이너스 회심공격시에 25%의 확률로 데미지가 2배가 된다."
$ws.Range("A8").Value = "비밀공작(裏稼業)"
$ws.Range("B8").Value = "암약(暗躍)"
$ws.Range("C8").Value = "10"
$ws.Range("D8").Value = "장전수UP(装填数UP) + 조합성공률+20%(調合成功率+20%) + 은밀(隠密)"
$ws.Range("A9").Value = "운신(運気)"
$ws.Range("B9").Value = "격운(激運)/강운(強運)/행운(幸運)/불운(不運)/재난(災難)"
$ws.Range("C9").Value = "20/15/10/-10/-15"
$ws.Range("D9").Value = "퀘스트 클리어 보수의 추첨에서 보수를 입수할 수 있는 확률이 통상의 22/32에서 31/32이 된다.`n/퀘스트 클리어 보수의 추첨에서 보수를 입수할 수 있는 확률이 통상의 22/32에서 28/32이 된다.`n/퀘스트 클리어 보수의 추첨에서 보수를 입수할 수 있는 확률이 통상의 22/32에서 25/32이 된다.`n/퀘스트 클리어 보수의 추첨에서 보수를 입수할 수 있는 확률이 통상의 22/32에서 16/32이 된다.`n/퀘스트 클리어 보수의 추첨에서 보수를 입수할 수 있는 확률이 통상의 22/32에서 8/32이 된다."
$ws.Range("A10").Value = "운반(運搬)"
$ws.Range("B10").Value = "운반의 달인(運搬の達人)"
$ws.Range("C10").Value = "10"
$ws.Range("D10").Value = "운반시의 이동속도가 상승. 높은 곳에서 뛰어 내려도 운반물을 잘 떨어뜨리지 않게 된다."
$ws.Range("A11").Value = "영웅의 방패(英雄の盾)"
$ws.Range("B11").Value = "영웅의 수호(英雄の護り)"
$ws.Range("C11").Value = "10"
$ws.Range("D11").Value = "작은 데미지(5이하)를 무효화. 상태이상, 아이템, 도트데미지 공격은 무효화 불가능"
$ws.Range("A12").Value = "SP연장(ＳＰ延長)"
$ws.Range("B12").Value = "SP시간연장(ＳＰ時間延長)"
$ws.Range("C12").Value = "10"
$ws.Range("D12").Value = "SP수기로 발동하는 SP상태의 효과시간이 1.25배가 된다."
$ws.Range("A13").Value = "염열적응(炎熱適応)"
$ws.Range("B13").Value = "남풍의 사냥꾼(南風の狩人)"
$ws.Range("C13").Value = "10"
$ws.Range("D13").Value = "더위 무효와 뜨거운 장소에서 공격력15와 방어력20상승. 핫드링크를 마시면 추가로 방어력10 상승"
$ws.Range("A14").Value = "오마(鏖魔)"
$ws.Range("B14").Value = "오마의 혼(鏖魔の魂)"
$ws.Range("C14").Value = "10"
$ws.Range("D14").Value = "공격력UP【大】(攻撃力UP【大】) + 회피성능+2(回避性能+2)"
$ws.Range("A15").Value = "대설주(大雪主)"
$ws.Range("B15").Value = "대설주의 혼(大雪主の魂)"
$ws.Range("C15").Value = "10"
$ws.Range("D15").Value = "회피거리UP(回避距離UP) + 스태미너 급속회복(スタミナ急速回復)"
$ws.Range("A16").Value = "농은(朧隠)"
$ws.Range("B16").Value = "농은의 혼(朧隠の魂)"
$ws.Range("C16").Value = "10"
$ws.Range("D16").Value = "고급 귀마개(高級耳栓) + 장전속도+3(装填速度+3) + 심검일체(心剣一体)"
$ws.Range("A17").Value = "가드강화(ガード強化)"
$ws.Range("B17").Value = "가드강화(ガード強化)"
$ws.Range("C17").Value = "10"
$ws.Range("D17").Value = "통상 가드 불가능한 공격이 가드 가능하게 된다. 특정 공격을 가드했을 때 일부 아이템 사용불가상태, 방어DOWN 상태가 되지 않는다."
$ws.Range("A18").Value = "가드성능(ガード性能)"
$ws.Range("B18").Value = "가드성능+2(ガード性能+2)/가드성능+1(ガード性能+1)/가드성능-1(ガード性能-1)"
$ws.Range("C18").Value = "15/10/-10"
$ws.Range("D18").Value = "가드 시 밀림 방지, 가드 시 위력 20경감/가드 시 밀림 방지, 가드 시 위력 10경감/가드 시 밀리기 쉬워짐, 가드 시 위력 20증가"
$ws.Range("A19").Value = "회심강화(会心強化)"
$ws.Range("B19").Value = "초회심(超会心)"
$ws.Range("C19").Value = "10"
$ws.Range("D19").Value = "회심 시 데미지가 통상의 1.25배에서 1.4배로 증가한다."
$ws.Range("A20").Value = "회피거리(回避距離)"
$ws.Range("B20").Value = "회피거리UP(回避距離UP)"
$ws.Range("C20").Value = "10"
$ws.Range("D20").Value = "구르기나 스텝의 이동거리가 1.5배로 늘어난다."
$ws.Range("A21").Value = "회피술(回避術)"
$ws.Range("B21").Value = "곡예사(軽業師)"
$ws.Range("C21").Value = "10"
$ws.Range("D21").Value = "체술+1(体術+1) + 회피성능+1(回避性能+1)"
$ws.Range("A22").Value = "회피성능(回避性能)"
$ws.Range("B22").Value = "회피성능+2(回避性能+2)/회피성능+1(回避性能+1)/회피성능DOWN(回避性能DOWN)"
$ws.Range("C22").Value = "15/10/-10"
$ws.Range("D22").Value = "구르기나 스텝의 무적시간이 통상의 0.2초에서 0.4초로 연장된다./구르기나 스텝의 무적시간이 통상의 0.2초에서 0.33초로 연장된다./구르기나 스텝의 무적시간이 통상의 0.2초에서 0.1초로 단축된다."
$ws.Range("A23").Value = "회복속도(回復速度)"
$ws.Range("B23").Value = "회복속도+2(回復速度+2)/회복속도+1(回復速度+1)/회복속도-1(回復速度-1)/회복속도-2(回復速度-2)"
$ws.Range("C23").Value = "15/10/-10/-15"
$ws.Range("D23").Value = "붉은 게이지의 회복속도가 4배가 된다/붉은 게이지의 회복속도가 2배가 된다/붉은 게이지의 회복속도가 1/2이 된다/붉은 게이지의 회복속도가 1/4이 된다"
$ws.Range("A24").Value = "회복량(回復量)"
$ws.Range("B24").Value = "체력회복량UP(体力回復量UP)/체력회복량DOWN(体力回復量DOWN)"
$ws.Range("C24").Value = "10/-10"
$ws.Range("D24").Value = "체력회복 아이템의 회복량이 1.25배로 늘어난다./체력회복 아이템의 회복량이 0.75배로 줄어든다."
$ws.Range("A25").Value = "확산탄추가(拡散弾追加)"
$ws.Range("B25").Value = "확산탄 전LV 추가(拡散弾全LV追加)/확산탄 LV1 추가(拡散弾LV1追加)"
$ws.Range("C25").Value = "15/10"
$ws.Range("D25").Value = "전LV의 확산탄을 사용할 수 있게 된다./LV1확산탄을 사용할 수 있게 된다."
$ws.Range("A26").Value = "가호(加護)"
$ws.Range("B26").Value = "정령의 가호(精霊の加護)/악령의 가호(悪霊の加護)"
$ws.Range("C26").Value = "10/-10"
$ws.Range("D26").Value = "1/4의 확률로 받는 데미지가 30% 감소한다./1/4의 확률로 받는 데미지가 30% 증가한다."
$ws.Range("A27").Value = "참기(我慢)"
$ws.Range("B27").Value = "호시탐탐(虎視眈々)"
$ws.Range("C27").Value = "10"
$ws.Range("D27").Value = "몬스터로부터 데미지를 받을 때 수기, 브레이브(데미지의 2배), 연금게이지가 쌓인다. 스트라이커의 경우 수기 게이지의 축적량이 데미지의 1.4배"
$ws.Range("A28").Value = "뇌속성 공격(雷属性攻撃)"
$ws.Range("B28").Value = "뇌속성 공격 강화+2(雷属性攻撃強化+2)/뇌속성 공격 강화+1(雷属性攻撃強化+1)/뇌속성 공격 약화(雷属性攻撃弱化)"
$ws.Range("C28").Value = "15/10/-10"
$ws.Range("D28").Value = "뇌속성의 공격이 1.1배+5이 된다./뇌속성의 공격이 1.05배+4가 된다./뇌속성의 공격이 0.75배가 된다."
$ws.Range("A29").Value = "뇌내성(雷耐性)"
$ws.Range("B29").Value = "뇌내성【대】(雷耐性【大】)/뇌내성【소】(雷耐性【小】)/뇌내성 약화(雷耐性弱化)"
$ws.Range("C29").Value = "15/10/-10"
$ws.Range("D29").Value = "雷耐性+20。合計耐性が25以上で雷属性やられ小と大を無効化/雷耐性+15。合計耐性が15以上で雷属性やられ小を無効化/雷耐性-20"
$ws.Range("A30").Value = "사냥꾼(狩人)"
$ws.Range("B30").Value = "헌터 생활"
$ws.Range("C30").Value = "10"
$ws.Range("D30").Value = "잘 구운 고기로 굽기 쉽게 된다. 지도가 없어도 맵이 표시된다. 낚시에서 물고기가 최초의 당김에  물게 된다. 낚시 미끼 아이템의 조합이 반드시 성공한다."

# Row heights for the two multi-line effect cells
$ws.Rows(6).RowHeight = 45
$ws.Rows(9).RowHeight = 75

# Wrap text on the long translated effect cells that now contain line breaks
$ws.Range("D6").WrapText = $true
$ws.Range("D9").WrapText = $true

# Column widths (approximate best-fit sizing matching the translated content)
$ws.Columns("A:A").ColumnWidth = 26.0833333333333
$ws.Columns("B:B").ColumnWidth = 106.416666666667
$ws.Columns("C:C").ColumnWidth = 20.75
$ws.Columns("D:D").ColumnWidth = 254.916666666667

# Restore selection / scroll position
$ws.Range("D29").Select()
